# issue #5: stock data output to json file
#
# Insert a new "property_category" column into the 股票 (stock) sheet,
# right after the "total" column and before the "date" column, with a
# constant value of "stock" for every data row. All the columns that used
# to live at/after that position (date, legislator_name, legislator_id)
# shift one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票 (stock) sheet

# Insert a brand-new column H, pushing the old H (date), I (legislator_name)
# and J (legislator_id) columns one to the right (-> I, J, K).
$ws.Range("H1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("H1").Value = "property_category"

# Value for the existing data row.
$ws.Range("H2").Value = "stock"
